$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.528.22"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "1.668.18"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06155"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06984"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.68"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.355"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5790"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.12"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "25.531.79"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006709"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").Value = "1.881.26"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.429"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.730"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.212"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.60"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.93"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.382"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.706"
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.09"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.006"
$ws.Range("E30").Value = "  +5.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07803"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.622"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04288"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.625"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9481"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6014"
$ws.Range("E36").Value = "  +4.55%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9353"
$ws.Range("E37").Value = "  +15.58%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.521"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.0000"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "102.00"
$ws.Range("E40").Value = "  +5.30%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01470"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.830"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3726"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.911"
$ws.Range("E44").Value = "  +4.64%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1110"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.154"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05262"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.68"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.375"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9990"
$ws.Range("E51").Value = "  +0.11%  "
